$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") moves from 2023-09-03 (45172) to 2023-09-06 (45175)
# for every existing data row (2..114).
for ($r = 2; $r -le 114; $r++) {
    $ws.Cells.Item($r, 3).Value = 45175
}

# Row 114 picks up an explicit 15pt row height.
$ws.Rows.Item(114).RowHeight = 15

# A new case row (115) is appended at the bottom of the table.
$row = 115

$ws.Cells.Item($row, 1).Value = "A 40983-2023"

$ws.Cells.Item($row, 2).Value = 45173
$ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 3).Value = 45175
$ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 4).Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Cells.Item($row, 5).Value = "ALE"
$ws.Cells.Item($row, 6).Value = "Kyrkan"

$ws.Cells.Item($row, 7).Value = 12.6
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0

$ws.Cells.Item($row, 18).WrapText = 1
$ws.Cells.Item($row, 18).Value = ""
